$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '''27.736.77'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '

$c = $ws.Range("D3")
$c.Value = '''1.847.81'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.04%  '

$ws.Range("E4").Value = '  -0.02%  '

$c = $ws.Range("D5")
$c.Value = '''313.93'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("E6").Value = '  +0.10%  '

$c = $ws.Range("D7")
$c.Value = '''0.4331'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.72%  '

$c = $ws.Range("D8")
$c.Value = '''0.3651'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.21%  '

$c = $ws.Range("D9")
$c.Value = '''45.02'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.92%  '

$c = $ws.Range("D10")
$c.Value = '''0.07335'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.87%  '

$c = $ws.Range("D11")
$c.Value = '''0.8762'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.77%  '

$c = $ws.Range("D12")
$c.Value = '''20.76'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.04%  '

$c = $ws.Range("D13")
$c.Value = '''1.825.54'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.73%  '

$c = $ws.Range("D14")
$c.Value = '''5.343'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("E15").Value = '  -0.84%  '

$c = $ws.Range("D16")
$c.Value = '''0.06922'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.06%  '

$c = $ws.Range("D17")
$c.Value = '''1.002'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.07%  '

$c = $ws.Range("D18")
$c.Value = '''79.95'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.39%  '

$c = $ws.Range("D19")
$c.Value = '''0.000008988'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.01%  '

$ws.Range("E20").Value = '  -0.07%  '

$c = $ws.Range("D21")
$c.Value = '''15.37'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.68%  '

$c = $ws.Range("D22")
$c.Value = '''27.604.39'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("E23").Value = '  +0.19%  '

$c = $ws.Range("D24")
$c.Value = '''10.39'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.17%  '

$c = $ws.Range("D25")
$c.Value = '''2.027.73'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.64%  '

$c = $ws.Range("D26")
$c.Value = '''1.981'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -3.11%  '

$c = $ws.Range("D27")
$c.Value = '''155.98'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.30%  '

$ws.Range("E28").Value = '  +1.81%  '

$c = $ws.Range("D29")
$c.Value = '''120.47'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +8.75%  '

$c = $ws.Range("D30")
$c.Value = '''5.247'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.66%  '

$c = $ws.Range("D31")
$c.Value = '''1.855'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.40%  '

$c = $ws.Range("D32")
$c.Value = '''0.08904'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.09%  '

$c = $ws.Range("D33")
$c.Value = '''0.7512'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.63%  '

$c = $ws.Range("D34")
$c.Value = '''4.538'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.42%  '

$c = $ws.Range("D35")
$c.Value = '''2.965'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.16%  '

$c = $ws.Range("D36")
$c.Value = '''1.120'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +3.06%  '

$ws.Range("E37").Value = '  +0.93%  '

$c = $ws.Range("D38")
$c.Value = '''0.05405'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.14%  '

$c = $ws.Range("D39")
$c.Value = '''0.01929'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.00%  '

$c = $ws.Range("D40")
$c.Value = '''2.845'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.57%  '

$ws.Range("E41").Value = '  +0.29%  '

$c = $ws.Range("D42")
$c.Value = '''0.1655'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.75%  '

$c = $ws.Range("D43")
$c.Value = '''6.656'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.42%  '

$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("E45").Value = '  +1.16%  '

$c = $ws.Range("D46")
$c.Value = '''0.06537'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.44%  '

$c = $ws.Range("D47")
$c.Value = '''0.4662'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.32%  '

$c = $ws.Range("D48")
$c.Value = '''104.31'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.20%  '

$c = $ws.Range("D49")
$c.Value = '''1.001'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.04%  '

$c = $ws.Range("D50")
$c.Value = '''1.623'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.10%  '

$c = $ws.Range("D51")
$c.Value = '''64.29'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.01%  '
